$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.31139146475743473,
    -0.0099999989425256786,
    -0.0089999989256561719,
    -0.011999999743572687,
    -0.0059999989404317944,
    -0.0059999989113954655,
    -0.019999998735920954,
    -0.019999998729548274,
    -0.0059999988970096396,
    -0.0059999988928964854,
    -0.0044999989110259264,
    -0.046714489894200906,
    -0.0059999989047945235,
    0.048572899846663908,
    -0.0059999989096759521,
    -0.0059999989098988848,
    -0.0059999989117827113,
    0.012578419747601544,
    -0.0089999989548896764,
    -0.0089999989458995344,
    -0.0089999989446525319,
    -0.008999998943906462,
    -0.0089999989172815376,
    -0.041999998494822499,
    -0.041999998487141532,
    0.027738696345132041,
    -0.0059999989042522905,
    -0.0059999988807817317,
    -0.011999998792299849,
    -0.019999998685406251,
    -0.02162232322986668,
    -0.020999998659728014,
    -0.0059999988449481734
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
